$p = $ppt.ActivePresentation

# --- Slide 1: Subtitle text tweak ("stop!" -> "stop !") ---
$slide1 = $p.Slides.Item(1)
$subtitleShape = $slide1.Shapes.Item(2)
$subtitleShape.TextFrame.TextRange.Text = "Once you start, you don’t want to stop !"

# --- Slide 4: "Technologies used" paragraph - trim trailing text and split run ---
$slide4 = $p.Slides.Item(4)
$contentShape = $slide4.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$techParagraph = $tr.Paragraphs(2, 1)

$firstPart = "	Nodejs, Express.JS, Sequelize, express-handlebars, MySQL, HTML, CSS, Bootswatch, JavaScript,                      	jQuery, Passport, Stripe "
$secondPart = "client API."

$techParagraph.Text = $firstPart
$techParagraph.InsertAfter($secondPart) | Out-Null
